$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "b"
$ws.Range("J2").Value = "Acknowledge (Backchannel)"
$ws.Range("I6").Value = "b"
$ws.Range("J6").Value = "Acknowledge (Backchannel)"
$ws.Range("I12").Value = "sd"
$ws.Range("J12").Value = "Statement-non-opinion"
$ws.Range("I13").Value = "sd"
$ws.Range("J13").Value = "Statement-non-opinion"
$ws.Range("I21").Value = "aa"
$ws.Range("J21").Value = "Agree/Accept"
$ws.Range("I26").Value = "sd"
$ws.Range("J26").Value = "Statement-non-opinion"
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"
$ws.Range("I40").Value = "ba"
$ws.Range("J40").Value = "Appreciation"
$ws.Range("I42").Value = "%"
$ws.Range("J42").Value = "Uninterpretable"
$ws.Range("I43").Value = "qy"
$ws.Range("J43").Value = "Yes-No-Question"
$ws.Range("I47").Value = "b"
$ws.Range("J47").Value = "Acknowledge (Backchannel)"
$ws.Range("I53").Value = "sd"
$ws.Range("J53").Value = "Statement-non-opinion"
$ws.Range("I62").Value = "sd"
$ws.Range("J62").Value = "Statement-non-opinion"
$ws.Range("I73").Value = "sv"
$ws.Range("J73").Value = "Statement-opinion"
$ws.Range("I74").Value = "qy"
$ws.Range("J74").Value = "Yes-No-Question"
$ws.Range("I88").Value = "aa"
$ws.Range("J88").Value = "Agree/Accept"
$ws.Range("I93").Value = "sd"
$ws.Range("J93").Value = "Statement-non-opinion"
$ws.Range("I95").Value = "sd"
$ws.Range("J95").Value = "Statement-non-opinion"
$ws.Range("I97").Value = "sd"
$ws.Range("J97").Value = "Statement-non-opinion"
$ws.Range("I98").Value = "%"
$ws.Range("J98").Value = "Uninterpretable"
$ws.Range("I99").Value = "sd"
$ws.Range("J99").Value = "Statement-non-opinion"
$ws.Range("I103").Value = "sd"
$ws.Range("J103").Value = "Statement-non-opinion"
$ws.Range("I111").Value = "sd"
$ws.Range("J111").Value = "Statement-non-opinion"
$ws.Range("I113").Value = "sd"
$ws.Range("J113").Value = "Statement-non-opinion"
$ws.Range("I117").Value = "sd"
$ws.Range("J117").Value = "Statement-non-opinion"
$ws.Range("I128").Value = "aa"
$ws.Range("J128").Value = "Agree/Accept"
$ws.Range("I150").Value = "sv"
$ws.Range("J150").Value = "Statement-opinion"
$ws.Range("I151").Value = "sd"
$ws.Range("J151").Value = "Statement-non-opinion"
$ws.Range("I160").Value = "sv"
$ws.Range("J160").Value = "Statement-opinion"
$ws.Range("I176").Value = "sd"
$ws.Range("J176").Value = "Statement-non-opinion"
$ws.Range("I181").Value = "sv"
$ws.Range("J181").Value = "Statement-opinion"
$ws.Range("I185").Value = "sd"
$ws.Range("J185").Value = "Statement-non-opinion"
$ws.Range("I186").Value = "aa"
$ws.Range("J186").Value = "Agree/Accept"
$ws.Range("I189").Value = "sv"
$ws.Range("J189").Value = "Statement-opinion"
$ws.Range("I200").Value = "ba"
$ws.Range("J200").Value = "Appreciation"
$ws.Range("I214").Value = "%"
$ws.Range("J214").Value = "Uninterpretable"
$ws.Range("I219").Value = "sv"
$ws.Range("J219").Value = "Statement-opinion"
$ws.Range("I221").Value = "qy"
$ws.Range("J221").Value = "Yes-No-Question"
$ws.Range("I225").Value = "sv"
$ws.Range("J225").Value = "Statement-opinion"
$ws.Range("I227").Value = "aa"
$ws.Range("J227").Value = "Agree/Accept"
$ws.Range("I233").Value = "ba"
$ws.Range("J233").Value = "Appreciation"
$ws.Range("I243").Value = "sd"
$ws.Range("J243").Value = "Statement-non-opinion"
$ws.Range("I263").Value = "b"
$ws.Range("J263").Value = "Acknowledge (Backchannel)"
$ws.Range("I266").Value = "ba"
$ws.Range("J266").Value = "Appreciation"
$ws.Range("I267").Value = "aa"
$ws.Range("J267").Value = "Agree/Accept"
$ws.Range("I272").Value = "sd"
$ws.Range("J272").Value = "Statement-non-opinion"
$ws.Range("I293").Value = "ba"
$ws.Range("J293").Value = "Appreciation"
$ws.Range("I295").Value = "sd"
$ws.Range("J295").Value = "Statement-non-opinion"
$ws.Range("I306").Value = "sd"
$ws.Range("J306").Value = "Statement-non-opinion"
$ws.Range("I318").Value = "aa"
$ws.Range("J318").Value = "Agree/Accept"
$ws.Range("I324").Value = "%"
$ws.Range("J324").Value = "Uninterpretable"
$ws.Range("I335").Value = "sd"
$ws.Range("J335").Value = "Statement-non-opinion"
$ws.Range("I346").Value = "aa"
$ws.Range("J346").Value = "Agree/Accept"
$ws.Range("I351").Value = "aa"
$ws.Range("J351").Value = "Agree/Accept"
$ws.Range("I352").Value = "sd"
$ws.Range("J352").Value = "Statement-non-opinion"
$ws.Range("I355").Value = "sd"
$ws.Range("J355").Value = "Statement-non-opinion"
$ws.Range("I360").Value = "sd"
$ws.Range("J360").Value = "Statement-non-opinion"
$ws.Range("I362").Value = "sd"
$ws.Range("J362").Value = "Statement-non-opinion"
$ws.Range("I366").Value = "aa"
$ws.Range("J366").Value = "Agree/Accept"
